$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New BGT voxel "LWP758" replaces the prior placeholder LWP754 tag in the header row
$ws.Range("BA1").Value = "LWP758"

# Row 23 holds the newly processed BGT/WM voxel data for LWP758, covering all 86
# metric columns (A:CH) across the three repeated LWP/Row/Col/Slice + metrics blocks.
$rowValues = @(
    "LWP758",
    "5",
    "5",
    "1",
    "7.647991e+00",
    "1.217171e+00",
    "7.237035e+00",
    "1.974221e+01",
    "4.281611e+00",
    "3.616100e+00",
    "0.000000e+00",
    "1.046911e+01",
    "1.222028e+01",
    "7.120622e+00",
    "0.000000e+00",
    "5.252108e+00",
    "1.549001e+01",
    "9.250815e+00",
    "6.931388e+00",
    "0.000000e+00",
    "1.525038e+00",
    "1.073994e+01",
    "3.621569e+00",
    "2.474082e+01",
    "1.054749e+01",
    "2.268940e+01",
    "LWP758",
    "5",
    "5",
    "1",
    "3.245176e+00",
    "2.087482e+00",
    "7.136637e+00",
    "2.210137e+00",
    "6.634073e+00",
    "3.206725e+00",
    "3.321025e+00",
    "7.975593e+00",
    "5.856856e+00",
    "3.586366e+00",
    "3.826889e+00",
    "3.122414e+00",
    "4.263797e+00",
    "5.489644e+00",
    "2.722585e+00",
    "9.076109e-01",
    "1.158232e+00",
    "6.309373e+00",
    "2.372040e+00",
    "2.264586e+00",
    "1.003864e+00",
    "8.736255e+00",
    "LWP758",
    "5",
    "5",
    "1",
    "1.207545e+00",
    "2.433514e+01",
    "9.590888e+01",
    "0.000000e+00",
    "4.202680e-02",
    "5.371094e+00",
    "7.020637e+00",
    "8.477733e+00",
    "8.477733e+00",
    "1.236412e-02",
    "1.194178e-02",
    "4.650955e+00",
    "1.673508e+01",
    "6.515555e+01",
    "3.347589e-02",
    "-3.779603e+01",
    "-9.397741e+00",
    "2.972641e+01",
    "7.579457e+00",
    "5.930642e-02",
    "1.464844e+00",
    "2.489279e+03",
    "1.845296e-01",
    "1.342773e+00",
    "1.872792e-02",
    "3.471696e-02",
    "-3.585251e-02",
    "1.632873e-03",
    "5.336861e+00",
    "4.519424e+00"
)

$rowNum = 23
for ($col = 1; $col -le $rowValues.Length; $col++) {
    $cell = $ws.Cells.Item($rowNum, $col)
    $cell.Value = "'" + $rowValues[$col - 1]
    $cell.Style = "Normal"
}
